# Replace the "word" column (B2:B193) values with the new working-set sequence.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'zögern',
    'sammeln',
    'lockern',
    'stillen',
    'freuen',
    'wüten',
    'dienen',
    'reizen',
    'weichen',
    'träumen',
    'streichen',
    'schauen',
    'bauen',
    'spüren',
    'werfen',
    'deuten',
    'graben',
    'holen',
    'meinen',
    'zünden',
    'runden',
    'schwingen',
    'kehren',
    'greifen',
    'bergen',
    'schwimmen',
    'jubeln',
    'proben',
    'kosten',
    'klettern',
    'stehlen',
    'trauen',
    'pflanzen',
    'lesen',
    'drohen',
    'sparen',
    'retten',
    'mühen',
    'sorgen',
    'dringen',
    'kichern',
    'schlucken',
    'wehtun',
    'krachen',
    'planen',
    'zeigen',
    'schenken',
    'sperren',
    'pressen',
    'saufen',
    'fallen',
    'flehen',
    'mauern',
    'fällen',
    'schwächen',
    'hassen',
    'brauchen',
    'schützen',
    'platzen',
    'stimmen',
    'kürzen',
    'knarren',
    'füttern',
    'drucken',
    'feiern',
    'trennen',
    'wahren',
    'löschen',
    'warnen',
    'wenden',
    'führen',
    'beißen',
    'geben',
    'streifen',
    'scheinen',
    'lügen',
    'tollen',
    'süßen',
    'liefern',
    'feuern',
    'hören',
    'backen',
    'erben',
    'passen',
    'siegen',
    'rücken',
    'quälen',
    'morden',
    'helfen',
    'rühren',
    'folgen',
    'gleichen',
    'handeln',
    'spinnen',
    'suchen',
    'heilen',
    'fließen',
    'herzen',
    'testen',
    'ändern',
    'jagen',
    'münzen',
    'fangen',
    'lenken',
    'betteln',
    'treiben',
    'wachsen',
    'wundern',
    'messen',
    'wirken',
    'sterben',
    'parken',
    'schrecken',
    'räumen',
    'töten',
    'spannen',
    'gelten',
    'baden',
    'grüßen',
    'fragen',
    'fahren',
    'altern',
    'boxen',
    'füllen',
    'seufzen',
    'beten',
    'stechen',
    'machen',
    'hacken',
    'achten',
    'formen',
    'stecken',
    'tropfen',
    'mögen',
    'leihen',
    'liegen',
    'bellen',
    'sinken',
    'sichern',
    'wecken',
    'malen',
    'flüchten',
    'schlagen',
    'decken',
    'irren',
    'werden',
    'biegen',
    'rufen',
    'schneiden',
    'rauben',
    'filmen',
    'albern',
    'lächeln',
    'bitten',
    'rasen',
    'gründen',
    'spielen',
    'heben',
    'weinen',
    'küssen',
    'fischen',
    'klingen',
    'schreiten',
    'starren',
    'bluten',
    'zielen',
    'drehen',
    'loben',
    'kriegen',
    'schwören',
    'heulen',
    'sprengen',
    'drücken',
    'stammen',
    'lieben',
    'hindern',
    'nennen',
    'ehren',
    'kümmern',
    'steuern',
    'enden',
    'schmecken',
    'hauen',
    'sprechen',
    'wetten',
    'arten',
    'scheitern',
    'ärgern',
    'kratzen',
    'äußern',
    'pfeifen',
    'schulden'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}

